# backlog sprint 3 wednesday
# Fill in the burndown numbers for Sprint 3 through Wednesday of week 2
# (columns N = Tue, O = Wed), and update the tail end of week 1
# (columns H..M) to reflect the day-by-day progress.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 3")

# --- Row 2 ("set up backlog") ---
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 6
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 5

# --- Row 3 ("connect to site") ---
$ws.Range("N3").Value = 2
$ws.Range("O3").Value = 2

# --- Row 4 ("code app") ---
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = 3

# --- Row 5 ("design UI") ---
$ws.Range("N5").Value = 2
$ws.Range("O5").Value = 1

# Row 11/12 hold SUM()/estimate formulas and recalc automatically.

# --- Update the active selection to reflect where the author left off ---
$ws.Activate()
$ws.Range("N4").Select()
